$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely; all rows below shift up by one.
$ws.Rows.Item(2).Delete()
